$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.368348956108093
$ws.Range("B1").Value = 1.938032746315002
$ws.Range("C1").Value = 2.8238525390625
$ws.Range("D1").Value = 4.861639022827148
$ws.Range("E1").Value = 1.020480036735535
